# Add new working-experience rows to the "Jobs" sheet of the CV workbook,
# fix a couple of small typos/spacing in the existing row, and leave the
# "Jobs" sheet as the active tab/selection (as it was when the edit was made).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jobs")

# --- Fix up the existing ERC Postdoc row (row 3) ---
$ws.Cells.Item(3, 3).Value = "2016 - Now"
$ws.Cells.Item(3, 5).Value = 'ERC Project: "Combinatorial Aspects of Computational Geometry"'

# --- New row 4: Web Administrator ---
# (B4 previously held an empty, underlined placeholder cell - clear that
#  leftover formatting now that it holds real text)
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "Web Administrator"
$ws.Cells.Item(4, 2).Font.Underline = $false
$ws.Cells.Item(4, 3).Value = "2013 - Now"
$ws.Cells.Item(4, 4).Value = "Olimpiada Mexicana de Matemáticas"
$ws.Cells.Item(4, 5).Value = "Server set-up, domain management and annual update of official website"

# --- New row 5: Chair ---
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "Chair"
$ws.Cells.Item(5, 3).Value = "2016 - 2018"
$ws.Cells.Item(5, 4).Value = "Asian Pacific Mathematical Olympiad"
$ws.Cells.Item(5, 5).Value = "Lead the organization of annual competition for students in 45 countries <br> Set up of official website to publish resullts and statistics."

# --- New row 6: Jury Member and Coordinator ---
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "Jury Member and Coordinator"
$ws.Cells.Item(6, 3).Value = "2012 - 2017"
$ws.Cells.Item(6, 4).Value = "International Mathematical Olympiad"
$ws.Cells.Item(6, 5).Value = "<ul> <li> 2012-2015 Member of the Jury of the International Mathematical Olympiad. Select problems and marking schemes for the competition. Grading and coordination of scores as Team Leader for Mexico </li> <li> 2017 Coordinator of the International Matematical Olympiad. Ellaboration of marking schemes and grading of international students </li> </ul>"

# --- New row 7: Audiovisual Translator ---
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "Audiovisual Translator"
$ws.Cells.Item(7, 3).Value = "2012 - 2014"
$ws.Cells.Item(7, 4).Value = "Khan Academy via Fundación Slim"
$ws.Cells.Item(7, 5).Value = "Translation from English to Spanish of 550 videos in mathematics, finance and physics"

# --- New row 8: Adjunct Professor and Teaching Assistant ---
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "Adjunct Professor and Teaching Assistant"
$ws.Cells.Item(8, 3).Value = "2010 - 2013"
$ws.Cells.Item(8, 4).Value = "Universidad Nacional Autónoma de México"
$ws.Cells.Item(8, 5).Value = "Courses: Analytic Geometry, Calculus, Complex Analysis, Probability, Problem-Solving Seminar, Real Analysis, Stochastic Processes"

# --- View/selection bookkeeping: leave the Programming sheet's selection at
#     H18, then finish with the Jobs sheet active and A9 selected (below the
#     newly-entered rows), mirroring the author's on-screen state at save time ---
$wsProg = $wb.Worksheets.Item("Programming")
$wsProg.Activate()
$wsProg.Range("H18").Select()

$ws.Activate()
$ws.Range("A9").Select()
